# Update the "剩余" (remaining days) and "开始时间" (start date) columns.
# For each data row (2..99):
#   - if the current remaining value (column E) is 1, it rolls over to 10
#     and the start date (column F) advances by 10 (e.g. 20251102 -> 20251112)
#   - otherwise the remaining value simply decreases by 1 and the start
#     date stays unchanged
# Row 36 contains a malformed date value and is left untouched, matching
# the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $eCell = $ws.Cells.Item($r, 5)
    $eValue = $eCell.Value2

    if ($null -eq $eValue) {
        continue
    }

    if ($eValue -eq 1) {
        $eCell.Value2 = 10

        $fCell = $ws.Cells.Item($r, 6)
        $fValue = $fCell.Value2
        if ($null -ne $fValue) {
            $fCell.Value2 = $fValue + 10
        }
    }
    else {
        $eCell.Value2 = $eValue - 1
    }
}
